# Fix Training Data Issue
# The "Date" column (BF) held a malformed value "4-23-2013-14" (an artifact
# of the filename "4-23-2013-14" being copied verbatim into the data rows).
# It should hold the correct ISO-style date string "2014-04-23" for the
# game date of 4/23 in the 2013-14 NBA season, stored as TEXT (not an
# Excel date serial) - exactly like the original malformed value was text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "4-23-2013-14"
$newValue = "2014-04-23"

# Data rows 2 through 31 all carry this value in column BF ("Date").
for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)  # column BF = 58
    if ($cell.Value() -eq $oldValue) {
        # Force text storage so Excel doesn't reinterpret the ISO-looking
        # string as a date serial number; then drop back to the sheet's
        # default (unstyled) cell style so no stray formatting lingers.
        $cell.NumberFormat = "@"
        $cell.Value = $newValue
        $cell.Style = "Normal"
    }
}
